$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.471.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "'1.895.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("E4").Value = "  -0.63%  "

$ws.Range("D5").Value = "'247.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.96%  "

$ws.Range("E6").Value = "  -4.53%  "

$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("D8").Value = "'44.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.84%  "

$ws.Range("D9").Value = "'0.352"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.09%  "

$ws.Range("E10").Value = "  -3.51%  "

$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("D12").Value = "'13.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("D13").Value = "'2.172.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("E14").Value = "  -0.89%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'1.911.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'4.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.30%  "

$ws.Range("D17").Value = "'35.470.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "'73.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.50%  "

$ws.Range("D19").Value = "'0.0₃0821"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.07%  "

$ws.Range("D20").Value = "'247.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("E21").Value = "  -2.37%  "

$ws.Range("E22").Value = "  -3.36%  "

$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("E24").Value = "  +5.76%  "

$ws.Range("E25").Value = "  -9.78%  "

$ws.Range("D26").Value = "'165.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("E27").Value = "  -2.63%  "

$ws.Range("E28").Value = "  -2.44%  "

$ws.Range("E29").Value = "  -4.39%  "

$ws.Range("D30").Value = "'4.128.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  +7.41%  "

$ws.Range("D32").Value = "'4.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.93%  "

$ws.Range("E33").Value = "  -1.94%  "

$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("E35").Value = "  -0.67%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.50%  "

$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -19.79%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.0680"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.72%  "

$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'17.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.65%  "

$ws.Range("D41").Value = "'97.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("E43").Value = "  -3.30%  "

$ws.Range("D44").Value = "'1.287.19"
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = "  -3.78%  "

$ws.Range("D46").Value = "'0.0799"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.96%  "

$ws.Range("E47").Value = "  -0.97%  "

$ws.Range("D48").Value = "'2.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.89%  "

$ws.Range("D49").Value = "'12.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.71%  "

$ws.Range("D50").Value = "'6.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.71%  "

$ws.Range("D51").Value = "'43.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.87%  "
